$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.664768333333334
$ws.Range("H2").Value = 19.994305
$ws.Range("I2").Value = 0.06516174319532789
$ws.Range("J2").Value = 0.0651617431953279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.135193666666667
$ws.Range("N2").Value = 3.405581
$ws.Range("O2").Value = 0.153770120695047
$ws.Range("P2").Value = 0.153770120695047
$ws.Range("Q2").Value = 7.565802801800556
$ws.Range("R2").Value = 68.09222521620499
$ws.Range("S2").Value = 0.01001992911584523
$ws.Range("T2").Value = 0.01001992911584523

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.664768333333334
$ws.Range("H3").Value = 19.994305
$ws.Range("I3").Value = 0.06516174319532789
$ws.Range("J3").Value = 0.0651617431953279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.287366
$ws.Range("N3").Value = 6.862098
$ws.Range("O3").Value = 0.3098401235152652
$ws.Range("P3").Value = 0.3098401235152652
$ws.Range("Q3").Value = 15.24476448354334
$ws.Range("R3").Value = 137.20288035189
$ws.Range("S3").Value = 0.02018972256011039
$ws.Range("T3").Value = 0.02018972256011039

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.664768333333334
$ws.Range("H4").Value = 19.994305
$ws.Range("I4").Value = 0.06516174319532789
$ws.Range("J4").Value = 0.0651617431953279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.138589666666667
$ws.Range("N4").Value = 9.415769000000001
$ws.Range("O4").Value = 0.425144471843918
$ws.Range("P4").Value = 0.425144471843918
$ws.Range("Q4").Value = 20.91797302172722
$ws.Range("R4").Value = 188.261757195545
$ws.Range("S4").Value = 0.02770315489520669
$ws.Range("T4").Value = 0.0277031548952067

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.664768333333334
$ws.Range("H5").Value = 19.994305
$ws.Range("I5").Value = 0.06516174319532789
$ws.Range("J5").Value = 0.0651617431953279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8212579999999999
$ws.Range("N5").Value = 2.463774
$ws.Range("O5").Value = 0.1112452839457698
$ws.Range("P5").Value = 0.1112452839457698
$ws.Range("Q5").Value = 5.473494311896666
$ws.Range("R5").Value = 49.26144880707
$ws.Range("S5").Value = 0.007248936624165585
$ws.Range("T5").Value = 0.007248936624165586

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.25592399999999
$ws.Range("H6").Value = 141.767772
$ws.Range("I6").Value = 0.4620233187619072
$ws.Range("J6").Value = 0.4620233187619072
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.135193666666667
$ws.Range("N6").Value = 3.405581
$ws.Range("O6").Value = 0.153770120695047
$ws.Range("P6").Value = 0.153770120695047
$ws.Range("Q6").Value = 53.64462563728132
$ws.Range("R6").Value = 482.8016307355319
$ws.Range("S6").Value = 0.07104538148994465
$ws.Range("T6").Value = 0.07104538148994463

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.25592399999999
$ws.Range("H7").Value = 141.767772
$ws.Range("I7").Value = 0.4620233187619072
$ws.Range("J7").Value = 0.4620233187619072
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.287366
$ws.Range("N7").Value = 6.862098
$ws.Range("O7").Value = 0.3098401235152652
$ws.Range("P7").Value = 0.3098401235152652
$ws.Range("Q7").Value = 108.091593856184
$ws.Range("R7").Value = 972.8243447056558
$ws.Range("S7").Value = 0.1431533621521221
$ws.Range("T7").Value = 0.1431533621521221

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.25592399999999
$ws.Range("H8").Value = 141.767772
$ws.Range("I8").Value = 0.4620233187619072
$ws.Range("J8").Value = 0.4620233187619072
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.138589666666667
$ws.Range("N8").Value = 9.415769000000001
$ws.Range("O8").Value = 0.425144471843918
$ws.Range("P8").Value = 0.425144471843918
$ws.Range("Q8").Value = 148.3169547551853
$ws.Range("R8").Value = 1334.852592796668
$ws.Range("S8").Value = 0.1964266598346052
$ws.Range("T8").Value = 0.1964266598346052

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.25592399999999
$ws.Range("H9").Value = 141.767772
$ws.Range("I9").Value = 0.4620233187619072
$ws.Range("J9").Value = 0.4620233187619072
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8212579999999999
$ws.Range("N9").Value = 2.463774
$ws.Range("O9").Value = 0.1112452839457698
$ws.Range("P9").Value = 0.1112452839457698
$ws.Range("Q9").Value = 38.80930563239199
$ws.Range("R9").Value = 349.2837506915279
$ws.Range("S9").Value = 0.05139791528523528
$ws.Range("T9").Value = 0.05139791528523528

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 20.98736333333333
$ws.Range("H10").Value = 62.96209
$ws.Range("I10").Value = 0.2051944060881897
$ws.Range("J10").Value = 0.2051944060881898
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.135193666666667
$ws.Range("N10").Value = 3.405581
$ws.Range("O10").Value = 0.153770120695047
$ws.Range("P10").Value = 0.153770120695047
$ws.Range("Q10").Value = 23.82472193603222
$ws.Range("R10").Value = 214.42249742429
$ws.Range("S10").Value = 0.03155276859012942
$ws.Range("T10").Value = 0.03155276859012942

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 20.98736333333333
$ws.Range("H11").Value = 62.96209
$ws.Range("I11").Value = 0.2051944060881897
$ws.Range("J11").Value = 0.2051944060881898
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.287366
$ws.Range("N11").Value = 6.862098
$ws.Range("O11").Value = 0.3098401235152652
$ws.Range("P11").Value = 0.3098401235152652
$ws.Range("Q11").Value = 48.00578131831334
$ws.Range("R11").Value = 432.05203186482
$ws.Range("S11").Value = 0.0635774601270062
$ws.Range("T11").Value = 0.06357746012700621

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 20.98736333333333
$ws.Range("H12").Value = 62.96209
$ws.Range("I12").Value = 0.2051944060881897
$ws.Range("J12").Value = 0.2051944060881898
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.138589666666667
$ws.Range("N12").Value = 9.415769000000001
$ws.Range("O12").Value = 0.425144471843918
$ws.Range("P12").Value = 0.425144471843918
$ws.Range("Q12").Value = 65.8707216885789
$ws.Range("R12").Value = 592.8364951972101
$ws.Range("S12").Value = 0.08723726740168984
$ws.Range("T12").Value = 0.08723726740168987

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 20.98736333333333
$ws.Range("H13").Value = 62.96209
$ws.Range("I13").Value = 0.2051944060881897
$ws.Range("J13").Value = 0.2051944060881898
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.8212579999999999
$ws.Range("N13").Value = 2.463774
$ws.Range("O13").Value = 0.1112452839457698
$ws.Range("P13").Value = 0.1112452839457698
$ws.Range("Q13").Value = 17.23604003640667
$ws.Range("R13").Value = 155.12436032766
$ws.Range("S13").Value = 0.02282690996936427
$ws.Range("T13").Value = 0.02282690996936427

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 27.37233166666667
$ws.Range("H14").Value = 82.116995
$ws.Range("I14").Value = 0.2676205319545753
$ws.Range("J14").Value = 0.2676205319545753
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.135193666666667
$ws.Range("N14").Value = 3.405581
$ws.Range("O14").Value = 0.153770120695047
$ws.Range("P14").Value = 0.153770120695047
$ws.Range("Q14").Value = 31.07289754989944
$ws.Range("R14").Value = 279.656077949095
$ws.Range("S14").Value = 0.04115204149912773
$ws.Range("T14").Value = 0.04115204149912772

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 27.37233166666667
$ws.Range("H15").Value = 82.116995
$ws.Range("I15").Value = 0.2676205319545753
$ws.Range("J15").Value = 0.2676205319545753
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.287366
$ws.Range("N15").Value = 6.862098
$ws.Range("O15").Value = 0.3098401235152652
$ws.Range("P15").Value = 0.3098401235152652
$ws.Range("Q15").Value = 62.61054079505667
$ws.Range("R15").Value = 563.49486715551
$ws.Range("S15").Value = 0.0829195786760266
$ws.Range("T15").Value = 0.0829195786760266

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 27.37233166666667
$ws.Range("H16").Value = 82.116995
$ws.Range("I16").Value = 0.2676205319545753
$ws.Range("J16").Value = 0.2676205319545753
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.138589666666667
$ws.Range("N16").Value = 9.415769000000001
$ws.Range("O16").Value = 0.425144471843918
$ws.Range("P16").Value = 0.425144471843918
$ws.Range("Q16").Value = 85.91051732157278
$ws.Range("R16").Value = 773.1946558941551
$ws.Range("S16").Value = 0.1137773897124163
$ws.Range("T16").Value = 0.1137773897124163

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 27.37233166666667
$ws.Range("H17").Value = 82.116995
$ws.Range("I17").Value = 0.2676205319545753
$ws.Range("J17").Value = 0.2676205319545753
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.8212579999999999
$ws.Range("N17").Value = 2.463774
$ws.Range("O17").Value = 0.1112452839457698
$ws.Range("P17").Value = 0.1112452839457698
$ws.Range("Q17").Value = 22.47974635990333
$ws.Range("R17").Value = 202.31771723913
$ws.Range("S17").Value = 0.02977152206700469
$ws.Range("T17").Value = 0.02977152206700469
